$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "team record" columns, matching the
# formatting already used by the other header cells in row 1 (AC1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 69   # AD
    $ws.Cells.Item($r, 31).Value = 93   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
